# Scheduled runner update: refresh market-price-derived columns (H-N)
# for the affected leve rows across the per-job-class profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 80.375
$ws.Range("I29").Value = 80.375
$ws.Range("K29").Value = 241.125
$ws.Range("M29").Value = 39.875

$ws.Range("H107").Value = 933.9355
$ws.Range("I107").Value = 969.75
$ws.Range("J107").Value = 599.6667
$ws.Range("K107").Value = 969.75
$ws.Range("L107").Value = 599.6667
$ws.Range("M107").Value = 950.25
$ws.Range("N107").Value = -4439.6667

$ws.Range("H113").Value = 3011.6667
$ws.Range("I113").Value = 2502.5
$ws.Range("J113").Value = 3157.1428
$ws.Range("K113").Value = 2502.5
$ws.Range("L113").Value = 3157.1428
$ws.Range("M113").Value = 751.5
$ws.Range("N113").Value = -9665.1428

$ws.Range("H132").Value = 1549.11
$ws.Range("I132").Value = 1358.3737
$ws.Range("J132").Value = 3477.6667
$ws.Range("K132").Value = 4075.1211
$ws.Range("L132").Value = 10433.0001
$ws.Range("M132").Value = -1545.1211
$ws.Range("N132").Value = -15493.0001

$ws.Range("H137").Value = 1320.7833
$ws.Range("I137").Value = 1131.641
$ws.Range("J137").Value = 1672.0476
$ws.Range("K137").Value = 3394.923
$ws.Range("L137").Value = 5016.142800000001
$ws.Range("M137").Value = -844.9230000000002
$ws.Range("N137").Value = -10116.1428

$ws.Range("H138").Value = 1449.56
$ws.Range("I138").Value = 721.5682
$ws.Range("J138").Value = 2021.5536
$ws.Range("K138").Value = 2164.7046
$ws.Range("L138").Value = 6064.6608
$ws.Range("M138").Value = 2975.2954
$ws.Range("N138").Value = -16344.6608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4193.28
$ws.Range("I32").Value = 3599.0652
$ws.Range("J32").Value = 11026.75
$ws.Range("K32").Value = 3599.0652
$ws.Range("L32").Value = 11026.75
$ws.Range("M32").Value = -3312.0652
$ws.Range("N32").Value = -11600.75

$ws.Range("H61").Value = 6291086
$ws.Range("I61").Value = 7408785.5
$ws.Range("J61").Value = 4025.25
$ws.Range("K61").Value = 7408785.5
$ws.Range("L61").Value = 4025.25
$ws.Range("M61").Value = -7408573.5
$ws.Range("N61").Value = -4449.25

$ws.Range("H110").Value = 27259.854
$ws.Range("I110").Value = 31837.414
$ws.Range("K110").Value = 31837.414
$ws.Range("M110").Value = -29792.414

$ws.Range("H122").Value = 46360.273
$ws.Range("I122").Value = 48425.047
$ws.Range("K122").Value = 145275.141
$ws.Range("M122").Value = -142825.141

$ws.Range("H136").Value = 6291086
$ws.Range("I136").Value = 7408785.5
$ws.Range("J136").Value = 4025.25
$ws.Range("K136").Value = 22226356.5
$ws.Range("L136").Value = 12075.75
$ws.Range("M136").Value = -22223806.5
$ws.Range("N136").Value = -17175.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5006
$ws.Range("I86").Value = 8326.666999999999
$ws.Range("J86").Value = 2515.5
$ws.Range("K86").Value = 8326.666999999999
$ws.Range("L86").Value = 2515.5
$ws.Range("M86").Value = -7203.666999999999
$ws.Range("N86").Value = -4761.5

$ws.Range("H89").Value = 5006
$ws.Range("I89").Value = 8326.666999999999
$ws.Range("J89").Value = 2515.5
$ws.Range("K89").Value = 41633.335
$ws.Range("L89").Value = 12577.5
$ws.Range("M89").Value = -36017.335
$ws.Range("N89").Value = -23809.5

$ws.Range("H105").Value = 10417896
$ws.Range("I105").Value = 12501268
$ws.Range("K105").Value = 12501268
$ws.Range("M105").Value = -12499521

$ws.Range("H134").Value = 1731.1014
$ws.Range("I134").Value = 1585.1111
$ws.Range("J134").Value = 2256.6667
$ws.Range("K134").Value = 4755.3333
$ws.Range("L134").Value = 6770.000100000001
$ws.Range("M134").Value = -2220.3333
$ws.Range("N134").Value = -11840.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1184.4117
$ws.Range("I58").Value = 838.3333
$ws.Range("J58").Value = 2015
$ws.Range("K58").Value = 838.3333
$ws.Range("L58").Value = 2015
$ws.Range("M58").Value = -635.3333
$ws.Range("N58").Value = -2421

$ws.Range("H134").Value = 3218.8215
$ws.Range("I134").Value = 3255.5908
$ws.Range("J134").Value = 3084
$ws.Range("K134").Value = 9766.7724
$ws.Range("L134").Value = 9252
$ws.Range("M134").Value = -7231.7724
$ws.Range("N134").Value = -14322

$ws.Range("H136").Value = 1184.4117
$ws.Range("I136").Value = 838.3333
$ws.Range("J136").Value = 2015
$ws.Range("K136").Value = 2514.9999
$ws.Range("L136").Value = 6045
$ws.Range("M136").Value = 35.0001000000002
$ws.Range("N136").Value = -11145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 462.2857
$ws.Range("I113").Value = 444.69565
$ws.Range("J113").Value = 472.4
$ws.Range("K113").Value = 1334.08695
$ws.Range("L113").Value = 1417.2
$ws.Range("M113").Value = 835.9130500000001
$ws.Range("N113").Value = -5757.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1846.4404
$ws.Range("I132").Value = 1524.8182
$ws.Range("J132").Value = 3025.7222
$ws.Range("K132").Value = 4574.4546
$ws.Range("L132").Value = 9077.1666
$ws.Range("M132").Value = -2044.4546
$ws.Range("N132").Value = -14137.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7336.3335
$ws.Range("I7").Value = 6002
$ws.Range("J7").Value = 10005
$ws.Range("K7").Value = 6002
$ws.Range("L7").Value = 10005
$ws.Range("M7").Value = -5890
$ws.Range("N7").Value = -10229

$ws.Range("H40").Value = 69449
$ws.Range("I40").Value = 113745
$ws.Range("J40").Value = 3005
$ws.Range("K40").Value = 113745
$ws.Range("L40").Value = 3005
$ws.Range("M40").Value = -113609
$ws.Range("N40").Value = -3277

$ws.Range("H108").Value = 29000
$ws.Range("J108").Value = 29000
$ws.Range("L108").Value = 29000
$ws.Range("N108").Value = -36680

$ws.Range("H126").Value = 7336.3335
$ws.Range("I126").Value = 6002
$ws.Range("J126").Value = 10005
$ws.Range("K126").Value = 18006
$ws.Range("L126").Value = 30015
$ws.Range("M126").Value = -15536
$ws.Range("N126").Value = -34955

$ws.Range("H132").Value = 3162.7917
$ws.Range("I132").Value = 2894.9429
$ws.Range("J132").Value = 3883.923
$ws.Range("K132").Value = 8684.8287
$ws.Range("L132").Value = 11651.769
$ws.Range("M132").Value = -6154.8287
$ws.Range("N132").Value = -16711.769

$ws.Range("H136").Value = 3877263.5
$ws.Range("I136").Value = 1090.5938
$ws.Range("J136").Value = 15153403
$ws.Range("K136").Value = 3271.7814
$ws.Range("L136").Value = 45460209
$ws.Range("M136").Value = -721.7814000000003
$ws.Range("N136").Value = -45465309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4052271
$ws.Range("I132").Value = 1229.1818
$ws.Range("J132").Value = 17158584
$ws.Range("K132").Value = 3687.5454
$ws.Range("L132").Value = 51475752
$ws.Range("M132").Value = -1157.5454
$ws.Range("N132").Value = -51480812

$ws.Range("H136").Value = 1648.7
$ws.Range("I136").Value = 1431.2877
$ws.Range("K136").Value = 4293.8631
$ws.Range("M136").Value = -1743.8631
